$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 4
$ws.Range("H4").Value = 1825.3334
$ws.Range("I4").Value = 240.5
$ws.Range("K4").Value = 240.5
$ws.Range("M4").Value = -126.5

# Row 19
$ws.Range("H19").Value = 2165.8572
$ws.Range("I19").Value = 705.75
$ws.Range("K19").Value = 705.75
$ws.Range("M19").Value = -530.75

# Row 132
$ws.Range("H132").Value = 3706.44
$ws.Range("J132").Value = 2115.3333
$ws.Range("L132").Value = 6345.999899999999
$ws.Range("N132").Value = -11405.9999

# Row 138
$ws.Range("H138").Value = 2201.172
$ws.Range("J138").Value = 2626.7834
$ws.Range("L138").Value = 7880.350199999999
$ws.Range("N138").Value = -18160.3502

$ws = $wb.Worksheets.Item("ARM")
# Row 24
$ws.Range("H24").Value = 69797.8
$ws.Range("J24").Value = 69797.8
$ws.Range("L24").Value = 69797.8
$ws.Range("N24").Value = -70545.8

# Row 32
$ws.Range("H32").Value = 27830098
$ws.Range("I32").Value = 33389648
$ws.Range("J32").Value = 32342.666
$ws.Range("K32").Value = 33389648
$ws.Range("L32").Value = 32342.666
$ws.Range("M32").Value = -33389361
$ws.Range("N32").Value = -32916.666

# Row 63
$ws.Range("H63").Value = 6953
$ws.Range("I63").Value = 7007.4287
$ws.Range("K63").Value = 7007.4287
$ws.Range("M63").Value = -6321.4287

# Row 66
$ws.Range("H66").Value = 6953
$ws.Range("I66").Value = 7007.4287
$ws.Range("K66").Value = 35037.14350000001
$ws.Range("M66").Value = -31605.14350000001

# Row 74
$ws.Range("H74").Value = 4820449.5
$ws.Range("I74").Value = 5816141.5
$ws.Range("J74").Value = 928197.9399999999
$ws.Range("K74").Value = 5816141.5
$ws.Range("L74").Value = 928197.9399999999
$ws.Range("M74").Value = -5815267.5
$ws.Range("N74").Value = -929945.9399999999

# Row 77
$ws.Range("H77").Value = 4820449.5
$ws.Range("I77").Value = 5816141.5
$ws.Range("J77").Value = 928197.9399999999
$ws.Range("K77").Value = 29080707.5
$ws.Range("L77").Value = 4640989.699999999
$ws.Range("M77").Value = -29076339.5
$ws.Range("N77").Value = -4649725.699999999

# Row 80
$ws.Range("H80").Value = 65005.5
$ws.Range("J80").Value = 65005.5
$ws.Range("L80").Value = 65005.5
$ws.Range("N80").Value = -67001.5

# Row 83
$ws.Range("H83").Value = 65005.5
$ws.Range("J83").Value = 65005.5
$ws.Range("L83").Value = 195016.5
$ws.Range("N83").Value = -205000.5

# Row 95
$ws.Range("H95").Value = 99908.75
$ws.Range("J95").Value = 99908.75
$ws.Range("L95").Value = 99908.75
$ws.Range("N95").Value = -105400.75

# Row 100
$ws.Range("H100").Value = 69797.8
$ws.Range("J100").Value = 69797.8
$ws.Range("L100").Value = 69797.8
$ws.Range("N100").Value = -71961.8

# Row 110
$ws.Range("H110").Value = 2856.8333
$ws.Range("I110").Value = 2856.8333
$ws.Range("K110").Value = 2856.8333
$ws.Range("M110").Value = -811.8332999999998

# Row 122
$ws.Range("H122").Value = 1825.625
$ws.Range("I122").Value = 1825.625
$ws.Range("K122").Value = 5476.875
$ws.Range("M122").Value = -3026.875

# Row 132
$ws.Range("H132").Value = 5629
$ws.Range("I132").Value = 2421.5
$ws.Range("K132").Value = 7264.5
$ws.Range("M132").Value = -4734.5

$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 1908.931
$ws.Range("J86").Value = 2668.6365
$ws.Range("L86").Value = 2668.6365
$ws.Range("N86").Value = -4914.636500000001

# Row 89
$ws.Range("H89").Value = 1908.931
$ws.Range("J89").Value = 2668.6365
$ws.Range("L89").Value = 13343.1825
$ws.Range("N89").Value = -24575.1825

# Row 92
$ws.Range("H92").Value = 67659.664
$ws.Range("J92").Value = 76490
$ws.Range("L92").Value = 76490
$ws.Range("N92").Value = -81482

$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Range("H7").Value = 213.59091
$ws.Range("I7").Value = 136.92857
$ws.Range("K7").Value = 136.92857
$ws.Range("M7").Value = -23.92857000000001

# Row 31
$ws.Range("H31").Value = 494595.66
$ws.Range("J31").Value = 1174119.9
$ws.Range("L31").Value = 1174119.9
$ws.Range("N31").Value = -1174709.9

# Row 34
$ws.Range("H34").Value = 494595.66
$ws.Range("J34").Value = 1174119.9
$ws.Range("L34").Value = 1174119.9
$ws.Range("N34").Value = -1174523.9

# Row 38
$ws.Range("H38").Value = 1371
$ws.Range("J38").Value = 1371
$ws.Range("L38").Value = 1371
$ws.Range("N38").Value = -2125

# Row 46
$ws.Range("H46").Value = 1371
$ws.Range("J46").Value = 1371
$ws.Range("L46").Value = 1371
$ws.Range("N46").Value = -1793

# Row 122
$ws.Range("H122").Value = 2101.6
$ws.Range("I122").Value = 2101.6
$ws.Range("K122").Value = 6304.799999999999
$ws.Range("M122").Value = -3854.799999999999

# Row 132
$ws.Range("H132").Value = 2837.4
$ws.Range("I132").Value = 2565.6843
$ws.Range("J132").Value = 8000
$ws.Range("K132").Value = 7697.0529
$ws.Range("L132").Value = 24000
$ws.Range("M132").Value = -5167.0529
$ws.Range("N132").Value = -29060

$ws = $wb.Worksheets.Item("CUL")
# Row 126
$ws.Range("H126").Value = 8227
$ws.Range("I126").Value = 7213.857
$ws.Range("K126").Value = 21641.571
$ws.Range("M126").Value = -16701.571

$ws = $wb.Worksheets.Item("GSM")
# Row 5
$ws.Range("H5").Value = 102309.625
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").ClearContents()

# Row 92
$ws.Range("H92").Value = 15871
$ws.Range("J92").Value = 15871
$ws.Range("L92").Value = 15871
$ws.Range("N92").Value = -19615

# Row 100
$ws.Range("H100").Value = 80659.664
$ws.Range("I100").Value = 79000
$ws.Range("K100").Value = 79000
$ws.Range("M100").Value = -77918

# Row 132
$ws.Range("H132").Value = 37043980
$ws.Range("I132").Value = 40003096
$ws.Range("K132").Value = 120009288
$ws.Range("M132").Value = -120006758

$ws = $wb.Worksheets.Item("LTW")
# Row 2
$ws.Range("H2").Value = 11033333
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 11033333
$ws.Range("K2").Value = 0
$ws.Range("M2").ClearContents()
$ws.Range("N2").Value = -11033557

# Row 7
$ws.Range("H7").Value = 40898.668
$ws.Range("I7").Value = 3377.4
$ws.Range("J7").Value = 148102.28
$ws.Range("K7").Value = 3377.4
$ws.Range("L7").Value = 148102.28
$ws.Range("M7").Value = -3265.4
$ws.Range("N7").Value = -148326.28

# Row 10
$ws.Range("H10").Value = 34400
$ws.Range("J10").Value = 34400
$ws.Range("L10").Value = 34400
$ws.Range("N10").Value = -34680

# Row 68
$ws.Range("H68").Value = 3270.5
$ws.Range("I68").Value = 2343.25
$ws.Range("K68").Value = 2343.25
$ws.Range("M68").Value = -1594.25

# Row 71
$ws.Range("H71").Value = 3270.5
$ws.Range("I71").Value = 2343.25
$ws.Range("K71").Value = 11716.25
$ws.Range("M71").Value = -7972.25

# Row 74
$ws.Range("H74").Value = 116000
$ws.Range("J74").Value = 116000
$ws.Range("L74").Value = 116000
$ws.Range("N74").Value = -117996

# Row 76
$ws.Range("H76").Value = 19444
$ws.Range("J76").Value = 19444
$ws.Range("L76").Value = 19444
$ws.Range("N76").Value = -20120

# Row 77
$ws.Range("H77").Value = 116000
$ws.Range("J77").Value = 116000
$ws.Range("L77").Value = 348000
$ws.Range("N77").Value = -357984

# Row 79
$ws.Range("H79").Value = 19444
$ws.Range("J79").Value = 19444
$ws.Range("L79").Value = 19444
$ws.Range("N79").Value = -21784

# Row 98
$ws.Range("H98").Value = 77754.5
$ws.Range("J98").Value = 77754.5
$ws.Range("L98").Value = 77754.5
$ws.Range("N98").Value = -83744.5

# Row 126
$ws.Range("H126").Value = 40898.668
$ws.Range("I126").Value = 3377.4
$ws.Range("J126").Value = 148102.28
$ws.Range("K126").Value = 10132.2
$ws.Range("L126").Value = 444306.84
$ws.Range("M126").Value = -7662.200000000001
$ws.Range("N126").Value = -449246.84

# Row 132
$ws.Range("H132").Value = 317675.9
$ws.Range("I132").Value = 5804.76
$ws.Range("J132").Value = 1431501.4
$ws.Range("K132").Value = 17414.28
$ws.Range("L132").Value = 4294504.199999999
$ws.Range("M132").Value = -14884.28
$ws.Range("N132").Value = -4299564.199999999

$ws = $wb.Worksheets.Item("WVR")
# Row 2
$ws.Range("H2").Value = 100833
$ws.Range("I2").Value = 117777.336
$ws.Range("J2").Value = 50000
$ws.Range("K2").Value = 117777.336
$ws.Range("L2").Value = 50000
$ws.Range("M2").Value = -117665.336
$ws.Range("N2").Value = -50224

# Row 95
$ws.Range("H95").Value = 92240.57000000001
$ws.Range("J95").Value = 92240.57000000001
$ws.Range("L95").Value = 92240.57000000001
$ws.Range("N95").Value = -97732.57000000001

# Row 97
$ws.Range("H97").Value = 89993
$ws.Range("J97").Value = 89993
$ws.Range("L97").Value = 89993
$ws.Range("N97").Value = -91975

# Row 115
$ws.Range("H115").Value = 161000
$ws.Range("J115").Value = 161000
$ws.Range("L115").Value = 161000
$ws.Range("N115").Value = -164134

# Row 132
$ws.Range("H132").Value = 2292.718
$ws.Range("I132").Value = 1494.2258
$ws.Range("J132").Value = 5386.875
$ws.Range("K132").Value = 4482.6774
$ws.Range("L132").Value = 16160.625
$ws.Range("M132").Value = -1952.6774
$ws.Range("N132").Value = -21220.625

# Row 136
$ws.Range("H136").Value = 1603.42
$ws.Range("I136").Value = 1166.425
$ws.Range("J136").Value = 3351.4
$ws.Range("K136").Value = 3499.275
$ws.Range("L136").Value = 10054.2
$ws.Range("M136").Value = -949.2749999999996
$ws.Range("N136").Value = -15154.2
